$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("AO2").Value = 22
$ws.Range("F2").Value = 2.68
$ws.Range("K2").Value = 3.85
$ws.Range("S2").Value = 2.98
$ws.Range("F6").Value = 2.28
$ws.Range("H6").Value = 2.9
$ws.Range("K6").Value = 500
$ws.Range("AE7").Value = 80
$ws.Range("AF7").Value = 10.5
$ws.Range("AJ7").Value = 15.5
$ws.Range("G7").Value = 1.6
$ws.Range("O7").Value = 1.23
$ws.Range("U7").Value = 2.1
$ws.Range("W7").Value = 2.66
$ws.Range("Q8").Value = 2.06
$ws.Range("AO9").Value = 190
$ws.Range("Q9").Value = 1.81
$ws.Range("U9").Value = 1.82
$ws.Range("AD10").Value = 230
$ws.Range("I10").Value = 17.5
$ws.Range("S10").Value = 2.08
$ws.Range("T10").Value = 1.99
$ws.Range("U10").Value = 1.88
$ws.Range("W10").Value = 5.2
$ws.Range("X10").Value = 95
$ws.Range("Y10").Value = 310
$ws.Range("AN11").Value = 8.4
$ws.Range("I11").Value = 4.8
$ws.Range("H12").Value = 6.2
$ws.Range("P12").Value = 2.32
$ws.Range("Q12").Value = 1.69
$ws.Range("S12").Value = 2.72
$ws.Range("AA13").Value = 220
$ws.Range("AD13").Value = 28
$ws.Range("AI13").Value = 90
$ws.Range("AL13").Value = 32
$ws.Range("AN13").Value = 6.8
$ws.Range("F13").Value = 1.51
$ws.Range("G13").Value = 1.53
$ws.Range("H13").Value = 7
$ws.Range("I13").Value = 7.8
$ws.Range("K13").Value = 5
$ws.Range("O13").Value = 1.23
$ws.Range("P13").Value = 2.34
$ws.Range("S13").Value = 2.74
$ws.Range("U13").Value = 2.12
$ws.Range("W13").Value = 2.88
$ws.Range("Y13").Value = 46
$ws.Range("Z13").Value = 65
$ws.Range("AA15").Value = 160
$ws.Range("AE15").Value = 1000
$ws.Range("AF15").Value = 12.5
$ws.Range("AH15").Value = 19.5
$ws.Range("AI15").Value = 1000
$ws.Range("AK15").Value = 17
$ws.Range("AM15").Value = 100
$ws.Range("AN15").Value = 7.2
$ws.Range("AO15").Value = 1000
$ws.Range("R15").Value = 1.56
$ws.Range("X15").Value = 22
$ws.Range("Z15").Value = 1000
$ws.Range("T16").Value = 1.63
$ws.Range("AG17").Value = 11
$ws.Range("T18").Value = 1.75
$ws.Range("AA20").Value = 980
$ws.Range("AC20").Value = 980
$ws.Range("AE20").Value = 980
$ws.Range("AG20").Value = 980
$ws.Range("AH20").Value = 980
$ws.Range("AJ20").Value = 150
$ws.Range("AL20").Value = 75
$ws.Range("AM20").Value = 140
$ws.Range("N20").Value = 3.55
$ws.Range("R20").Value = 1.34
$ws.Range("S20").Value = 3.3
$ws.Range("T20").Value = 1.81
$ws.Range("U20").Value = 1.97
$ws.Range("X20").Value = 17.5
$ws.Range("Y20").Value = 980
$ws.Range("Q21").Value = 1.48
$ws.Range("AG22").Value = 10
$ws.Range("J22").Value = 3.85
$ws.Range("H24").Value = 2.94
$ws.Range("AE25").Value = 30
$ws.Range("AK25").Value = 48
$ws.Range("AM25").Value = 100
$ws.Range("AN25").Value = 44
$ws.Range("AO25").Value = 20
$ws.Range("I26").Value = 10.5
$ws.Range("R26").Value = 1.51
$ws.Range("S26").Value = 2.8
$ws.Range("J27").Value = 4.1
$ws.Range("R27").Value = 1.57
$ws.Range("AG29").Value = 9.6
$ws.Range("T29").Value = 1.9
$ws.Range("AG30").Value = 9.6
$ws.Range("AH30").Value = 24
$ws.Range("R30").Value = 1.49
$ws.Range("T30").Value = 1.89
$ws.Range("N31").Value = 4.3
$ws.Range("AD32").Value = 9.6
$ws.Range("M32").Value = 1.04
$ws.Range("X32").Value = 25
$ws.Range("AE33").Value = 330
$ws.Range("AN33").Value = 3.9
$ws.Range("I33").Value = 18
$ws.Range("L33").Value = 1.27
$ws.Range("P33").Value = 2.6
$ws.Range("T33").Value = 2.3
$ws.Range("U33").Value = 1.7
$ws.Range("AJ34").Value = 16.5
$ws.Range("AN34").Value = 8
$ws.Range("U34").Value = 2.04
$ws.Range("W34").Value = 2.66
$ws.Range("T36").Value = 1.69
$ws.Range("P37").Value = 1.84
$ws.Range("AH38").Value = 980
$ws.Range("I38").Value = 2.08
$ws.Range("J38").Value = 3.6
$ws.Range("N38").Value = 3.75
$ws.Range("V38").Value = 1.93
$ws.Range("AC40").Value = 980
$ws.Range("AD40").Value = 980
$ws.Range("AF40").Value = 70
$ws.Range("AK40").Value = 140
$ws.Range("AL40").Value = 110
$ws.Range("AM40").Value = 160
$ws.Range("AO40").Value = 1000
$ws.Range("N40").Value = 4.5
$ws.Range("O40").Value = 1.23
$ws.Range("R40").Value = 1.48
$ws.Range("S40").Value = 2.7
$ws.Range("T40").Value = 1.89
$ws.Range("U40").Value = 1.94
$ws.Range("Y40").Value = 980
$ws.Range("Z40").Value = 980
$ws.Range("AG41").Value = 10
$ws.Range("G41").Value = 1.77
$ws.Range("W41").Value = 2.28
$ws.Range("G42").Value = 2.96
$ws.Range("H42").Value = 2.76
$ws.Range("I42").Value = 3.1
$ws.Range("J42").Value = 3.2
$ws.Range("W42").Value = 1.51
$ws.Range("G43").Value = 2.06
$ws.Range("H43").Value = 1.09
$ws.Range("I43").Value = 7.8
$ws.Range("J43").Value = 2.6
$ws.Range("V43").Value = 1.18
$ws.Range("W43").Value = 1.94
$ws.Range("F44").Value = 1.04
$ws.Range("F45").Value = 1.8
$ws.Range("J45").Value = 3.4
$ws.Range("N45").Value = 2.54
$ws.Range("T45").Value = 1.04
$ws.Range("W45").Value = 2.1
